$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "BY ADMINISTRASI" (row 8) and "BIAYA ADM KARTU" (row 9) detail
# rows; this shifts the trailing SALDO AKHIR summary row up from row 10 to
# row 8 and shrinks the used range to A1:H8.
$ws.Rows(8).Delete()
$ws.Rows(8).Delete()

# Helper: set a date-looking cell as plain text (not auto-converted to a
# date serial) and strip the resulting number-format style so the cell
# keeps its original (unstyled) appearance.
function Set-TextValue($addr, $value) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - SALDO AWAL
$ws.Range("G2").Value = 148943002

# Row 3
Set-TextValue "A3" "04-Jun-2024"
$ws.Range("C3").Value = "050   Tanpa Kategori"
$ws.Range("F3").Value = 114574866
$ws.Range("G3").Value = 263517868

# Row 4
Set-TextValue "A4" "30-Jun-2024"
$ws.Range("B4").Value = "JASA GIROIBUNGA"
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 183543
$ws.Range("G4").Value = 263701411

# Row 5
Set-TextValue "A5" "30-Jun-2024"
$ws.Range("B5").Value = "PPH"
$ws.Range("E5").Value = 36709
$ws.Range("G5").Value = 263664702

# Row 6
Set-TextValue "A6" "30-Jun-2024"
$ws.Range("B6").Value = "BY ADMINISTRASI"
$ws.Range("E6").Value = 12000
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 263652702

# Row 7
Set-TextValue "A7" "30-Jun-2024"
$ws.Range("B7").Value = "BIAYA ADM KARTU"
$ws.Range("E7").Value = 10000
$ws.Range("G7").Value = 263642702

# Row 8 - SALDO AKHIR (was row 10, now shifted to row 8 by the deletions
# above; clear the transaction fields and keep the running balance/label)
$ws.Range("A8").Value = ""
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 263642702
$ws.Range("H8").Value = "SALDO AKHIR"
